$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: add a header label in column A (row 1) and give it the same
# formatting as the existing header cells (bold / bordered / centered style
# that is already applied to B1:E1), without introducing a new style entry.
# ---------------------------------------------------------------------------
function Set-HeaderA1 {
    param($ws, [string]$text)
    $ws.Range("A1").Value = $text
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Helper: strip the bold/bordered header style from a label cell so it goes
# back to the default "Normal" style (removes the s="1" attribute).
# ---------------------------------------------------------------------------
function Clear-LabelStyle {
    param($ws, [string]$cellAddr)
    $ws.Range($cellAddr).Style = "Normal"
}

# ---------------------------------------------------------------------------
# Helper: set a cell's value to a *text* string that looks like a number
# (e.g. a year like "2015") without Excel auto-converting it to a numeric
# cell and without introducing a stray new style record. We build the text
# via a formula returning a string, then paste-special the computed value
# back over itself (values only), which collapses it to a plain shared
# string while preserving the cell's existing style.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($ws, [string]$cellAddr, [string]$text)
    $ws.Range($cellAddr).Formula = "=""$text"""
    $ws.Range($cellAddr).Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Sheets 1-4 ("Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
# "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)") all share the
# same row layout (rows 2-12) in column A. Add a header in A1 and fix the
# accentuation / formatting of the source labels.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    Set-HeaderA1 $ws "Fonte/Tecnologia"

    Clear-LabelStyle $ws "A2"

    $ws.Range("A3").Value = "Gás Natural"
    Clear-LabelStyle $ws "A3"

    $ws.Range("A4").Value = "Carvão"
    Clear-LabelStyle $ws "A4"

    Clear-LabelStyle $ws "A5"

    $ws.Range("A6").Value = "Óleos Comb"
    Clear-LabelStyle $ws "A6"

    Clear-LabelStyle $ws "A7"

    $ws.Range("A8").Value = "Eólica"
    Clear-LabelStyle $ws "A8"

    Clear-LabelStyle $ws "A9"

    Clear-LabelStyle $ws "A10"

    $ws.Range("A11").Value = "Pot. Compl."
    Clear-LabelStyle $ws "A11"

    Clear-LabelStyle $ws "A12"
}

# ---------------------------------------------------------------------------
# Sheet 5 - "Emissoes Totais (MtCO2eq)"
# Add header, fix labels, and remove the "Teto" row entirely.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

Set-HeaderA1 $ws5 "Período"

$ws5.Range("A2").Value = "P.Médio"
Clear-LabelStyle $ws5 "A2"

$ws5.Range("A3").Value = "P.Crítico"
Clear-LabelStyle $ws5 "A3"

$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# Sheet 6 - "Custo Total (bilhões de R$)"
# Add header, rename the cost column to "2015", fix labels and update values.
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

Set-HeaderA1 $ws6 "Tipo Expansão"

Set-TextValue $ws6 "B1" "2015"

$ws6.Range("A2").Value = "Expansão Centralizada"
Clear-LabelStyle $ws6 "A2"
$ws6.Range("B2").Value = 588

$ws6.Range("A3").Value = "Expansão por GD"
Clear-LabelStyle $ws6 "A3"
$ws6.Range("B3").Value = 99
